$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 360
$ws.Range("J2").Value = 1471
$ws.Range("K2").Value = 8
$ws.Range("L2").Value = 380
$ws.Range("M2").Value = 20
$ws.Range("N2").Value = 278
$ws.Range("Q2").Value = 3
$ws.Range("R2").Value = 15
$ws.Range("S2").Value = 158
$ws.Range("T2").Value = 277
$ws.Range("U2").Value = 18
$ws.Range("V2").Value = 2220
$ws.Range("X2").Value = 2294
$ws.Range("Y2").Value = 2
$ws.Range("Z2").Value = 42
$ws.Range("AA2").Value = 20
